$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    10 = @('sd', 'Statement-non-opinion')
    14 = @('b', 'Acknowledge (Backchannel)')
    26 = @('b', 'Acknowledge (Backchannel)')
    36 = @('sv', 'Statement-opinion')
    44 = @('b', 'Acknowledge (Backchannel)')
    46 = @('aa', 'Agree/Accept')
    49 = @('%', 'Uninterpretable')
    56 = @('aa', 'Agree/Accept')
    62 = @('sd', 'Statement-non-opinion')
    63 = @('sd', 'Statement-non-opinion')
    64 = @('sv', 'Statement-opinion')
    66 = @('%', 'Uninterpretable')
    67 = @('%', 'Uninterpretable')
    68 = @('%', 'Uninterpretable')
    73 = @('sd', 'Statement-non-opinion')
    85 = @('aa', 'Agree/Accept')
    93 = @('sv', 'Statement-opinion')
    103 = @('aa', 'Agree/Accept')
    104 = @('sd', 'Statement-non-opinion')
    106 = @('aa', 'Agree/Accept')
    112 = @('sd', 'Statement-non-opinion')
    115 = @('aa', 'Agree/Accept')
    116 = @('aa', 'Agree/Accept')
    132 = @('sd', 'Statement-non-opinion')
    137 = @('sd', 'Statement-non-opinion')
    138 = @('b', 'Acknowledge (Backchannel)')
    139 = @('sd', 'Statement-non-opinion')
    161 = @('sd', 'Statement-non-opinion')
    164 = @('b', 'Acknowledge (Backchannel)')
    168 = @('b', 'Acknowledge (Backchannel)')
    176 = @('ba', 'Appreciation')
    181 = @('sd', 'Statement-non-opinion')
    190 = @('sd', 'Statement-non-opinion')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
